$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 27-30: status changed from "Ej påbörjat" to "Klart"
$ws.Range("C27").Value2 = "Klart"
$ws.Range("C28").Value2 = "Klart"
$ws.Range("C29").Value2 = "Klart"
$ws.Range("C30").Value2 = "Klart"

# Rows 40-42: status changed from "Ej påbörjat" to "Klart"
$ws.Range("C40").Value2 = "Klart"
$ws.Range("C41").Value2 = "Klart"
$ws.Range("C42").Value2 = "Klart"

# Verklig tid (actual time) added for rows 41 and 42
$ws.Range("E41").Value2 = 5
$ws.Range("E42").Value2 = 2

# Row 50: status changed from "Ej påbörjat" to "Påbörjat", and Verklig tid set
$ws.Range("C50").Value2 = "Påbörjat"
$ws.Range("E50").Value2 = 4

# Update the selection to reflect the new active cell
$ws.Range("F51").Select()
